$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/medication-knowledge-cost-type"
$ws.Range("B3").Value = "8.0.0"
$ws.Range("B4").Value = "MedicationKnowledgeCostTypeCodeSystem"
$ws.Range("B5").Value = "Medication Knowledge Cost Type Codes"
$ws.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$ws.Range("B9").Value = "LinuxForHealth Team"
$ws.Range("B11").Value = "Medication knowledge cost type code"
